# Weekly update: insert a new data row above the current row 260,
# pushing the existing rows 260-275 down to 261-276.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 260 (shifts 260..275 down to 261..276)
$ws.Rows.Item(260).Insert()

# Populate the new row 260 with the new weekly record
$ws.Cells.Item(260, 1).Value2 = 10
$ws.Cells.Item(260, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(260, 3).Value2 = "La Araucanía"
$ws.Cells.Item(260, 4).Value2 = 44610
$ws.Cells.Item(260, 5).Value2 = 9
$ws.Cells.Item(260, 6).Value2 = 100112044
$ws.Cells.Item(260, 7).Value2 = "Perejil"
$ws.Cells.Item(260, 8).Value2 = "Sin especificar"
$ws.Cells.Item(260, 9).Value2 = "Primera"
$ws.Cells.Item(260, 10).Value2 = 30
$ws.Cells.Item(260, 11).Value2 = 4000
$ws.Cells.Item(260, 12).Value2 = 4000
$ws.Cells.Item(260, 13).Value2 = 4000
$ws.Cells.Item(260, 14).Value2 = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(260, 15).Value2 = "Provincia de Cautín"
$ws.Cells.Item(260, 16).Value2 = 1333
$ws.Cells.Item(260, 17).Value2 = 3
$ws.Cells.Item(260, 18).Value2 = "Hortaliza"
